# Add a "Slovakia" market tab, cloned from the existing "Portugal" tab,
# with the market name / ticket number swapped and the duplicate
# P32AR / P32DR repeater rows removed (Slovakia only lists PR1DS / PR8AS).

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Clone the Portugal sheet and place the copy right after it.
$portugal.Copy($null, $portugal)
$slovakia = $wb.Worksheets.Item("Portugal (2)")
$slovakia.Name = "Slovakia"

# Swap in the Slovakia-specific market name / ticket reference.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3177"

# Slovakia's repeater list skips P32AR / P32DR (rows 16-17), so remove
# them and let PR1DS / PR8AS shift up.
$slovakia.Rows("16:17").Delete()

# Rows 3-4 no longer need the taller wrapped-text height once the cell
# above was re-entered; restore their auto-fit height.
$slovakia.Rows("3:4").AutoFit()

# New sheet becomes the active/selected tab, with A8:A17 selected.
$slovakia.Select()
$slovakia.Range("A8:A17").Select()

# Portugal is no longer the active tab; its selection becomes "select all".
$portugal.Select()
$portugal.Cells.Select()

# Leave Slovakia as the active sheet (matches activeTab pointing at it).
$slovakia.Select()
